# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Rule row "R40" (cell B11 on sheet "Rules") is relabelled to "1".
# The new label must stay a text value (matching how the row's other
# rule labels - R10/R20/R30 - are stored as shared strings) rather than
# being auto-coerced into a number, and the cell's existing style must
# be preserved.
#
# A plain `Range.Value = "1"` assignment gets auto-typed as a numeric 1
# by the COM layer, so instead we enter a formula that evaluates to the
# text string "1", then convert that formula to its static value via a
# copy / paste-special(values) round-trip in place. That leaves B11
# holding the literal text "1" with its original number format/style
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")
$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)   # xlPasteValues
